$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 44.615383
$ws.Range("I11").Value = 44.615383
$ws.Range("K11").Value = 44.615383
$ws.Range("M11").Value = 95.38461699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1053.409
$ws.Range("J29").Value = 994.1177
$ws.Range("L29").Value = 2982.3531
$ws.Range("N29").Value = -3544.3531

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 383.78946
$ws.Range("I33").Value = 475.08334
$ws.Range("J33").Value = 227.28572
$ws.Range("K33").Value = 475.08334
$ws.Range("L33").Value = 227.28572
$ws.Range("M33").Value = -246.08334
$ws.Range("N33").Value = -685.28572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 626
$ws.Range("I39").Value = 30.25
$ws.Range("J39").Value = 3009
$ws.Range("K39").Value = 90.75
$ws.Range("L39").Value = 9027
$ws.Range("M39").Value = 205.25
$ws.Range("N39").Value = -9619

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9913.333000000001
$ws.Range("I51").Value = 9813.333000000001
$ws.Range("J51").Value = 9933.333000000001
$ws.Range("K51").Value = 9813.333000000001
$ws.Range("L51").Value = 9933.333000000001
$ws.Range("M51").Value = -9329.333000000001
$ws.Range("N51").Value = -10901.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4333
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4333
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4849.7
$ws.Range("I74").Value = 4356.857
$ws.Range("J74").Value = 5999.6665
$ws.Range("K74").Value = 4356.857
$ws.Range("L74").Value = 5999.6665
$ws.Range("M74").Value = -3420.857
$ws.Range("N74").Value = -7871.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2673.8572
$ws.Range("I76").Value = 2413
$ws.Range("K76").Value = 2413
$ws.Range("M76").Value = -2098

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4849.7
$ws.Range("I77").Value = 4356.857
$ws.Range("J77").Value = 5999.6665
$ws.Range("K77").Value = 21784.285
$ws.Range("L77").Value = 29998.3325
$ws.Range("M77").Value = -17104.285
$ws.Range("N77").Value = -39358.3325

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2673.8572
$ws.Range("I79").Value = 2413
$ws.Range("K79").Value = 2413
$ws.Range("M79").Value = -1321

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 942.5909
$ws.Range("I98").Value = 1025.8422
$ws.Range("K98").Value = 1025.8422
$ws.Range("M98").Value = 472.1578

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 942.5909
$ws.Range("I122").Value = 1025.8422
$ws.Range("K122").Value = 3077.5266
$ws.Range("M122").Value = -627.5266000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 899.3
$ws.Range("I131").Value = 665.8889
$ws.Range("K131").Value = 1997.6667
$ws.Range("M131").Value = 3042.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5760.04
$ws.Range("I132").Value = 1925.85
$ws.Range("K132").Value = 5777.549999999999
$ws.Range("M132").Value = -3247.549999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 4861.4
$ws.Range("I135").Value = 1671.7142
$ws.Range("J135").Value = 12304
$ws.Range("K135").Value = 15045.4278
$ws.Range("L135").Value = 110736
$ws.Range("M135").Value = -12510.4278
$ws.Range("N135").Value = -115806

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 20199522
$ws.Range("I137").Value = 836502.8
$ws.Range("K137").Value = 2509508.4
$ws.Range("M137").Value = -2506958.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5430.5884
$ws.Range("I138").Value = 1474.4348
$ws.Range("J138").Value = 6898.1934
$ws.Range("K138").Value = 4423.3044
$ws.Range("L138").Value = 20694.5802
$ws.Range("M138").Value = 716.6956
$ws.Range("N138").Value = -30974.5802

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14800.184
$ws.Range("I32").Value = 14028.648
$ws.Range("J32").Value = 21744
$ws.Range("K32").Value = 14028.648
$ws.Range("L32").Value = 21744
$ws.Range("M32").Value = -13741.648
$ws.Range("N32").Value = -22318

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4546.8184
$ws.Range("I61").Value = 4741.5
$ws.Range("K61").Value = 4741.5
$ws.Range("M61").Value = -4529.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1799.6976
$ws.Range("I132").Value = 1709.7
$ws.Range("K132").Value = 5129.1
$ws.Range("M132").Value = -2599.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4546.8184
$ws.Range("I136").Value = 4741.5
$ws.Range("K136").Value = 14224.5
$ws.Range("M136").Value = -11674.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 900.4
$ws.Range("I22").Value = 778
$ws.Range("K22").Value = 778
$ws.Range("M22").Value = -605

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3821.182
$ws.Range("I105").Value = 3770.9473
$ws.Range("J105").Value = 4139.3335
$ws.Range("K105").Value = 3770.9473
$ws.Range("L105").Value = 4139.3335
$ws.Range("M105").Value = -2023.9473
$ws.Range("N105").Value = -7633.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1560.6364
$ws.Range("I22").Value = 1461.3334
$ws.Range("J22").Value = 1679.8
$ws.Range("K22").Value = 1461.3334
$ws.Range("L22").Value = 1679.8
$ws.Range("M22").Value = -1111.3334
$ws.Range("N22").Value = -2379.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 45676
$ws.Range("J52").Value = 48160
$ws.Range("L52").Value = 48160
$ws.Range("N52").Value = -48748

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2572
$ws.Range("J94").Value = 3210.9167
$ws.Range("L94").Value = 3210.9167
$ws.Range("N94").Value = -4112.9167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 59999
$ws.Range("I104").Value = 59999
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 59999
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -57378
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H115").Value = 49248
$ws.Range("J115").Value = 49248
$ws.Range("L115").Value = 49248
$ws.Range("N115").Value = -51598

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 97589.60000000001
$ws.Range("J139").Value = 97362.5
$ws.Range("L139").Value = 97362.5
$ws.Range("N139").Value = -107642.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7000104.5
$ws.Range("J4").Value = 34500300
$ws.Range("L4").Value = 103500900
$ws.Range("N4").Value = -103501124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 3199.4285
$ws.Range("J92").Value = 3199.4285
$ws.Range("L92").Value = 3199.4285
$ws.Range("N92").Value = -6943.4285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 39498
$ws.Range("J96").Value = 39498
$ws.Range("L96").Value = 39498
$ws.Range("N96").Value = -44990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 90000
$ws.Range("J114").Value = 90000
$ws.Range("L114").Value = 90000
$ws.Range("N114").Value = -98678

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 289284.38
$ws.Range("I122").Value = 436346.2
$ws.Range("K122").Value = 1309038.6
$ws.Range("M122").Value = -1306588.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 145858.42
$ws.Range("I132").Value = 251164.75
$ws.Range("K132").Value = 753494.25
$ws.Range("M132").Value = -750964.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 90181.336
$ws.Range("J133").Value = 90181.336
$ws.Range("L133").Value = 90181.336
$ws.Range("N133").Value = -100301.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13482923
$ws.Range("I40").Value = 4809348
$ws.Range("K40").Value = 4809348
$ws.Range("M40").Value = -4809212

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2808.8948
$ws.Range("I82").Value = 1679.75
$ws.Range("J82").Value = 3630.0908
$ws.Range("K82").Value = 1679.75
$ws.Range("L82").Value = 3630.0908
$ws.Range("M82").Value = -1318.75
$ws.Range("N82").Value = -4352.0908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2808.8948
$ws.Range("I85").Value = 1679.75
$ws.Range("J85").Value = 3630.0908
$ws.Range("K85").Value = 1679.75
$ws.Range("L85").Value = 3630.0908
$ws.Range("M85").Value = -431.75
$ws.Range("N85").Value = -6126.0908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7113.528
$ws.Range("I122").Value = 4088.6667
$ws.Range("K122").Value = 12266.0001
$ws.Range("M122").Value = -9816.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3810.602
$ws.Range("I132").Value = 3022.0364
$ws.Range("K132").Value = 9066.109199999999
$ws.Range("M132").Value = -6536.109199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 29257
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 353.86957
$ws.Range("I113").Value = 430.3846
$ws.Range("J113").Value = 254.4
$ws.Range("K113").Value = 1291.1538
$ws.Range("L113").Value = 763.2
$ws.Range("M113").Value = 878.8462
$ws.Range("N113").Value = -5103.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3123.625
$ws.Range("I126").Value = 3024.5
$ws.Range("K126").Value = 9073.5
$ws.Range("M126").Value = -6603.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 27779158
$ws.Range("I132").Value = 37038536
$ws.Range("K132").Value = 111115608
$ws.Range("M132").Value = -111113078

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 98394
$ws.Range("I138").Value = 98394
$ws.Range("K138").Value = 98394
$ws.Range("M138").Value = -93254
